# Station4 Q file: add cross-sectional Area alongside existing Discharge (Q)
# calculations, plus a small Atotal/Qtotal summary block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column headers -------------------------------------------------
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"

# --- D column (segment midpoint depth) was filled individually; redo it
#     as one fill-down so it becomes a shared formula across D3:D9, same
#     as the author's edit.
$ws.Range("D3:D9").Formula = "=(A3/100+(A4/100-A3/100)/2)"

# --- Area (G) column: first segment uses 0 as the left edge, the rest
#     are filled down from G4 (so G4:G15 becomes one shared formula,
#     matching G3 typed separately first).
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

# --- Atotal (H2) sums the Area column -----------------------------------
$ws.Range("H2").Formula = "=SUM(G2:G11)"

# --- Small summary block pulling the totals together --------------------
$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

# --- Match the saved selection on the new summary cells ------------------
$ws.Range("J2:K2").Select() | Out-Null
